$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: torsion spring ---
$ws.Range("E4").Value2 = "https://www.mcmaster.com/9271K403/"
$ws.Range("A4").Value2 = "1.25"" leg, 0.4"" OD, 0.25"" ID_Torsion Spring"
$ws.Range("B4").Value2 = 2
$ws.Range("C4").Value2 = 5.38
$ws.Range("C4").Font.Name = "Arial"
$ws.Range("C4").Font.Size = 9
$ws.Range("C4").Font.Color = 3355443

# --- Row 5: shoulder screw ---
$ws.Range("A5").Value2 = "1/4"" Shoulder Diameter, 2-1/2"" Shoulder Length, 10-24 Thread, shoulder screw"
$ws.Range("E5").Value2 = "https://www.mcmaster.com/91259A105/"
$ws.Range("B5").Value2 = 4
$ws.Range("C5").Value2 = 2.42
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Color = 3368499
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Rows("5").RowHeight = 42.75

# --- Row 6: 8-32 3/4" button head screw ---
$ws.Range("A6").Value2 = "8-32, 3/4"" Long button head screw"
$ws.Range("A7").Value2 = "8-32, 2"" Long button head screw"
$ws.Range("A8").Value2 = "8-32 thin nylock nuts"
$ws.Range("E6").Value2 = "https://www.mcmaster.com/92949A197/"
$ws.Range("E8").Value2 = "https://www.mcmaster.com/90633A009/"
$ws.Range("E7").Value2 = "https://www.mcmaster.com/92949A207/"

$ws.Range("B6").Value2 = 1
$ws.Range("C6").Value2 = 6.53
$ws.Range("B7").Value2 = 1
$ws.Range("C7").Value2 = 5.77
$ws.Range("B8").Value2 = 1
$ws.Range("C8").Value2 = 3.23

$ws.Range("A6").Font.Name = "Arial"
$ws.Range("A6").Font.Color = 3368499
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $true

$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Color = 3368499
$ws.Range("A7").VerticalAlignment = -4108
$ws.Range("A7").WrapText = $true

$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Color = 3368499
$ws.Range("A8").VerticalAlignment = -4108
$ws.Range("A8").WrapText = $true

# --- D12 becomes a plain (non-shared) formula ---
$ws.Range("D12").Formula = "=B12*C12"

# --- Totals column ---
$ws.Range("F1").Value2 = "total:"
$ws.Range("G1").Formula = "=SUM(C2:C19)"

# --- Selection / view ---
$ws.Range("F3").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
